# Import / Export XL Fixes #25
#
# Adds a new "ID" column as the first column on every worksheet, shifting
# all existing columns one place to the right, and populates the new
# column with the row's numeric id (header "ID" in row 1).

$wb = $excel.ActiveWorkbook

# --- Aclass --------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Aclass")
$ws1.Columns("A:A").Insert()
$ws1.Range("A1").Value = "ID"
$ws1.Range("A2").Value = 1
$ws1.Range("A3").Value = 3

# --- Bclass ----------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Bclass")
$ws2.Columns("A:A").Insert()
$ws2.Range("A1").Value = "ID"
$ws2.Range("A2").Value = 1
$ws2.Range("A3").Value = 2

# --- Dclass ------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Dclass")
$ws3.Columns("A:A").Insert()
$ws3.Range("A1").Value = "ID"
